$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "GRT-USD"
$ws.Range("A21").Value = "BSCX-USD"
